$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M (13th column). This shifts the old
# M:Q columns (DESCTO, SUBTOTAL, IVA, ICE, TOTAL) one column to the right
# (to N:R), updates the dimension/merge/row spans automatically, and adds
# blank cells in the new column M for rows 4-7.
$ws.Columns("M:M").Insert()

# Give the brand-new column M the same width as column L so it keeps an
# explicit width instead of falling back to the sheet default.
$ws.Columns("M:M").ColumnWidth = $ws.Columns("L:L").ColumnWidth()

# Center-align the new (currently blank) filter cells in the new column,
# matching the rest of the filter row cells.
$ws.Range("M4:M7").HorizontalAlignment = -4108

# Fill in the new column's header (row 9): "SUBTOTAL" now lives in M,
# swapping places with "DESCTO" which the column insert already moved to N.
$ws.Range("M9").Value = "SUBTOTAL"

# Insert the brand-new "TOT SIN IVA" header into column O, overwriting the
# "SUBTOTAL" text that the column insert shifted there from the old N9.
$ws.Range("O9").Value = "TOT SIN IVA"

# Match the saved selection from the edited workbook.
$ws.Range("A9").Select()
